# Bilan-dataSet.xlsx edit: add AI "Model prediction" column (H), move the
# "Calcule avec RMS" / "Déviation" computed columns from I/J to L/M, and
# append a duplicated block of G2:G14 values into F30:F42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New column H: "Model prediction" header + literal values (rows 2-14)
# ---------------------------------------------------------------------
$ws.Range("H1").Value = "Model prediction"

$modelPredictions = @{
    2  = 94.687248093435002
    3  = 94.959077458293095
    4  = 64.4323129265115
    5  = 94.016962524455096
    6  = 87.279647474527806
    7  = 85.653102188874399
    8  = 82.359686356429293
    9  = 85.242444694750503
    10 = 89.632383875928298
    11 = 89.348543617324793
    12 = 89.566795231854002
    13 = 89.590290789469705
    14 = 94.132884696235607
}
foreach ($r in 2..14) {
    $ws.Cells.Item($r, 8).Value = $modelPredictions[$r]
}

# Widen column H like the rest of the data columns (target stored width
# 19.5546875; the engine quantizes ColumnWidth to 1/6-character steps, so
# 18.7 is the closest input that lands on the nearest representable width)
$ws.Columns.Item(8).ColumnWidth = 18.7

# ---------------------------------------------------------------------
# 2) Move the existing "Calcule avec RMS" (I) / "Déviation" (J) columns
#    to L / M, rewriting every formula explicitly (no more shared
#    formulas) so the row-by-row values stay identical.
# ---------------------------------------------------------------------
$ws.Range("I1:J23").ClearContents()

$ws.Range("L1").Value = "Calcule avec RMS"
$ws.Range("M1").Value = "Déviation"

foreach ($r in 2..14) {
    $ws.Cells.Item($r, 12).Formula = "=110-25*((E$r/C$r)/(F$r/D$r))"
    $ws.Cells.Item($r, 13).Formula = "=ABS(L$r-G$r)"
}

$ws.Range("L23").Value = "Average Deviation"
$ws.Range("M23").Formula = "=AVERAGE(M2:M14)"

# ---------------------------------------------------------------------
# 3) New block of rows 30-42: F column re-lists the G2:G14 SaO2 values.
# ---------------------------------------------------------------------
$gValues = @{
    30 = 95.987200000000001
    31 = 94.021199999999993
    32 = 64.846500000000006
    33 = 93.720200000000006
    34 = 89.3994
    35 = 82.846500000000006
    36 = 81.807699999999997
    37 = 85.099699999999999
    38 = 91.972399999999993
    39 = 92.323599999999999
    40 = 89.424300000000002
    41 = 85.177300000000002
    42 = 94.146799999999999
}
foreach ($r in 30..42) {
    $ws.Cells.Item($r, 6).Value = $gValues[$r]
}

# ---------------------------------------------------------------------
# 4) Selection cosmetics (matches the authored sheetView selection)
# ---------------------------------------------------------------------
$ws.Range("I9").Select() | Out-Null
